# Apply Excalibur_Profits.xlsx market-data refresh across all class sheets
# (values sourced from scheduled market-board snapshot)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 30404.334
$ws.Range("J62").Value = 15606.5
$ws.Range("L62").Value = 15606.5
$ws.Range("N62").Value = -16854.5
$ws.Range("H65").Value = 30404.334
$ws.Range("J65").Value = 15606.5
$ws.Range("L65").Value = 78032.5
$ws.Range("N65").Value = -84272.5
$ws.Range("H70").Value = 4332.25
$ws.Range("I70").Value = 1620
$ws.Range("J70").Value = 6269.5713
$ws.Range("K70").Value = 4860
$ws.Range("L70").Value = 18808.7139
$ws.Range("M70").Value = -4590
$ws.Range("N70").Value = -19348.7139
$ws.Range("H73").Value = 4332.25
$ws.Range("I73").Value = 1620
$ws.Range("J73").Value = 6269.5713
$ws.Range("K73").Value = 4860
$ws.Range("L73").Value = 18808.7139
$ws.Range("M73").Value = -3924
$ws.Range("N73").Value = -20680.7139
$ws.Range("H96").Value = 2023.5555
$ws.Range("I96").Value = 1280.6666
$ws.Range("K96").Value = 3841.9998
$ws.Range("M96").Value = -2468.9998
$ws.Range("H116").Value = 265202.16
$ws.Range("I116").Value = 12347.5
$ws.Range("K116").Value = 12347.5
$ws.Range("M116").Value = -8905.5
$ws.Range("H125").Value = 5295.6
$ws.Range("J125").Value = 11455.5
$ws.Range("L125").Value = 103099.5
$ws.Range("N125").Value = -108019.5
$ws.Range("H137").Value = 1203378.1
$ws.Range("I137").Value = 1050.9412
$ws.Range("J137").Value = 3247334.2
$ws.Range("K137").Value = 3152.8236
$ws.Range("L137").Value = 9742002.600000001
$ws.Range("M137").Value = -602.8235999999997
$ws.Range("N137").Value = -9747102.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8599.317999999999
$ws.Range("J2").Value = 23858.143
$ws.Range("L2").Value = 23858.143
$ws.Range("N2").Value = -24084.143
$ws.Range("H31").Value = 7686
$ws.Range("I31").Value = 7686
$ws.Range("K31").Value = 7686
$ws.Range("M31").Value = -7392
$ws.Range("H32").Value = 5851069
$ws.Range("I32").Value = 6292434.5
$ws.Range("K32").Value = 6292434.5
$ws.Range("M32").Value = -6292147.5
$ws.Range("H45").Value = 2254.4614
$ws.Range("I45").Value = 2264.4546
$ws.Range("K45").Value = 2264.4546
$ws.Range("M45").Value = -1887.4546
$ws.Range("H61").Value = 1451090.2
$ws.Range("I61").Value = 1588860.8
$ws.Range("K61").Value = 1588860.8
$ws.Range("M61").Value = -1588648.8
$ws.Range("H110").Value = 676.3158
$ws.Range("I110").Value = 658.3333
$ws.Range("K110").Value = 658.3333
$ws.Range("M110").Value = 1386.6667
$ws.Range("H116").Value = 8599.317999999999
$ws.Range("J116").Value = 23858.143
$ws.Range("L116").Value = 23858.143
$ws.Range("N116").Value = -28446.143
$ws.Range("H122").Value = 2590.6072
$ws.Range("I122").Value = 2519.889
$ws.Range("K122").Value = 7559.667
$ws.Range("M122").Value = -5109.667
$ws.Range("H136").Value = 1451090.2
$ws.Range("I136").Value = 1588860.8
$ws.Range("K136").Value = 4766582.4
$ws.Range("M136").Value = -4764032.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8599.317999999999
$ws.Range("J3").Value = 23858.143
$ws.Range("L3").Value = 23858.143
$ws.Range("N3").Value = -24086.143
$ws.Range("H94").Value = 1034.3158
$ws.Range("I94").Value = 950.3077
$ws.Range("J94").Value = 1216.3334
$ws.Range("K94").Value = 950.3077
$ws.Range("L94").Value = 1216.3334
$ws.Range("M94").Value = -499.3077
$ws.Range("N94").Value = -2118.3334
$ws.Range("H105").Value = 2666.6667
$ws.Range("I105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("M105").Value = -253
$ws.Range("H134").Value = 960746
$ws.Range("I134").Value = 1273496.4
$ws.Range("K134").Value = 3820489.2
$ws.Range("M134").Value = -3817954.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 5643.3335
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 5643.3335
$ws.Range("K20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").Value = 16930.0005
$ws.Range("N20").Value = -17384.0005
$ws.Range("H41").Value = 194
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H129").Value = 1343.3334
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1343.3334
$ws.Range("K129").Value = 0
$ws.Range("L129").ClearContents()
$ws.Range("M129").Value = 4030.0002
$ws.Range("N129").Value = -14030.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6782.727
$ws.Range("I70").Value = 6955.8887
$ws.Range("J70").Value = 6003.5
$ws.Range("K70").Value = 6955.8887
$ws.Range("L70").Value = 6003.5
$ws.Range("M70").Value = -6685.8887
$ws.Range("N70").Value = -6543.5
$ws.Range("H73").Value = 6782.727
$ws.Range("I73").Value = 6955.8887
$ws.Range("J73").Value = 6003.5
$ws.Range("K73").Value = 6955.8887
$ws.Range("L73").Value = 6003.5
$ws.Range("M73").Value = -6019.8887
$ws.Range("N73").Value = -7875.5
$ws.Range("H97").Value = 5015.2856
$ws.Range("J97").Value = 4633.3335
$ws.Range("L97").Value = 4633.3335
$ws.Range("N97").Value = -5625.3335
$ws.Range("H102").Value = 3083.5715
$ws.Range("I102").Value = 2581.75
$ws.Range("K102").Value = 2581.75
$ws.Range("M102").Value = -959.75
$ws.Range("H122").Value = 27537.906
$ws.Range("J122").Value = 6093.524
$ws.Range("L122").Value = 18280.572
$ws.Range("N122").Value = -23180.572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 363.36365
$ws.Range("I22").Value = 299.7
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 299.7
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -4.699999999999989
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 363.36365
$ws.Range("I27").Value = 299.7
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 299.7
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -192.7
$ws.Range("N27").Value = -1214
$ws.Range("H40").Value = 3259.4736
$ws.Range("I40").Value = 2871.6155
$ws.Range("K40").Value = 2871.6155
$ws.Range("M40").Value = -2735.6155
$ws.Range("H46").Value = 2030.6923
$ws.Range("J46").Value = 2197.9092
$ws.Range("L46").Value = 2197.9092
$ws.Range("N46").Value = -2573.9092
$ws.Range("H68").Value = 7080.6
$ws.Range("I68").Value = 11500
$ws.Range("J68").Value = 4134.3335
$ws.Range("K68").Value = 11500
$ws.Range("L68").Value = 4134.3335
$ws.Range("M68").Value = -10751
$ws.Range("N68").Value = -5632.3335
$ws.Range("H71").Value = 7080.6
$ws.Range("I71").Value = 11500
$ws.Range("J71").Value = 4134.3335
$ws.Range("K71").Value = 57500
$ws.Range("L71").Value = 20671.6675
$ws.Range("M71").Value = -53756
$ws.Range("N71").Value = -28159.6675
$ws.Range("H136").Value = 51871.92
$ws.Range("I136").Value = 1506.9445
$ws.Range("J136").Value = 181381.86
$ws.Range("K136").Value = 4520.833500000001
$ws.Range("L136").Value = 544145.58
$ws.Range("M136").Value = -1970.833500000001
$ws.Range("N136").Value = -549245.58

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 182880
$ws.Range("I62").Value = 302133.34
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 302133.34
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -301509.34
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 182880
$ws.Range("I65").Value = 302133.34
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 1510666.7
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -1507546.7
$ws.Range("N65").Value = -26240
$ws.Range("H81").Value = 4089.9092
$ws.Range("I81").Value = 4127.2856
$ws.Range("K81").Value = 8254.5712
$ws.Range("M81").Value = -7193.5712
$ws.Range("H84").Value = 4089.9092
$ws.Range("I84").Value = 4127.2856
$ws.Range("K84").Value = 41272.856
$ws.Range("M84").Value = -35968.856
$ws.Range("H113").Value = 1632.125
$ws.Range("I113").Value = 346.6842
$ws.Range("K113").Value = 1040.0526
$ws.Range("M113").Value = 1129.9474
$ws.Range("H122").Value = 2842.6562
$ws.Range("I122").Value = 2497.7083
$ws.Range("J122").Value = 3877.5
$ws.Range("K122").Value = 7493.124899999999
$ws.Range("L122").Value = 11632.5
$ws.Range("M122").Value = -5043.124899999999
$ws.Range("N122").Value = -16532.5
